$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New tweets scraped for #BelajarDariRumah (rows 259-273) ---
$ws.Range("A259").Value = [double]"1.342467637910114E+18"
$ws.Range("B259").Value = "Kondisi pandemi Covid-19 ini membuat kegiatan belajar sekolah terpaksa harus dilakukan online di rumah, padahal masih banyak guru &amp; siswa yg belum memiliki gawai untuk belajar online. Untuk itulah program Garuda diluncurkan... #bantuan #BelajarDariRumah https://t.co/LUqG2V5nQ5"
$ws.Range("C259").Value = "AyoGerakBareng"
$ws.Range("D259").Value = "Fri Dec 25 13:50:01 +0000 2020"
$ws.Range("A260").Value = [double]"1.342259705075266E+18"
$ws.Range("B260").Value = "Selamat pagi, #SahabatDikbud. Beragam tayangan #BelajardariRumah di @TVRINasional siap hadir menjadi pilihan untuk menemani #SahabatDikbud di hari libur ini. Yuk, simak jadwal untuk hari Jumat, 25 Desember 2020! `n#MerdekaBelajar`n#BersamaHadapiKorona https://t.co/lcP4j4rYex"
$ws.Range("C260").Value = "Kemdikbud_RI"
$ws.Range("D260").Value = "Fri Dec 25 00:03:46 +0000 2020"
$ws.Range("A261").Value = [double]"1.3422078081297979E+18"
$ws.Range("B261").Value = "Jadwal Program ""BELAJAR DARI RUMAH"" Hari Jumat, 25 Desember 2020 pukul 08.00-11.00 WIB. Dilanjutkan dengan  Main-Main Serius pukul 10.30 WIB`n#SemangatBaruMarta`n#BelajarDariRumah`n#MediaPemersatuBangsa`n#TVRI #TVRINASIONAL https://t.co/uFo9xyOnzo"
$ws.Range("C261").Value = "TheMartaSaputra"
$ws.Range("D261").Value = "Thu Dec 24 20:37:33 +0000 2020"
$ws.Range("A262").Value = [double]"1.342163513901343E+18"
$ws.Range("B262").Value = "Jadwal Program ""BELAJAR DARI RUMAH"" Hari Jumat, 25 Desember 2020 pukul 08.00-11.00 WIB. Dilanjutkan dengan  Main-Main Serius pukul 10.30 WIB`n#BelajarDariRumah`n#MediaPemersatuBangsa`n#TVRI #TVRINASIONAL https://t.co/nsJI5jSXcM"
$ws.Range("C262").Value = "TVRINasional"
$ws.Range("D262").Value = "Thu Dec 24 17:41:32 +0000 2020"
$ws.Range("A263").Value = [double]"1.3420848396215749E+18"
$ws.Range("B263").Value = "Meskipun sekarang sudah mulai liburan, jangan sampai kita lupa hari ya #GenPrestasi !`nTebak-tebakan yuk! Ada yang tahu jawaban dari pertanyaan yang ada di gambar?`n#IndiHomeStudyByIndiHome #BelajarDariRumah #dirumahaja #KamisKuis #MenghitungHari https://t.co/J736qxj8pP"
$ws.Range("C263").Value = "indihome_study"
$ws.Range("D263").Value = "Thu Dec 24 12:28:55 +0000 2020"
$ws.Range("A264").Value = [double]"1.3418891005849521E+18"
$ws.Range("B264").Value = "Selamat pagi, #SahabatDikbud. Yuk, isi liburanmu dengan aktivitas yang menyenangkan #dirumahsaja! Salah satunya, dengan menonton beragam tayangan menarik #BelajardariRumah di @TVRINasional. Berikut jadwal acaranya. `n#MerdekaBelajar`n#BersamaHadapiKorona https://t.co/lueuLrFOUS"
$ws.Range("C264").Value = "Kemdikbud_RI"
$ws.Range("D264").Value = "Wed Dec 23 23:31:07 +0000 2020"
$ws.Range("A265").Value = [double]"1.3418689778577779E+18"
$ws.Range("B265").Value = "Jadwal Program ""BELAJAR DARI RUMAH""  Hari Kamis, 24 Desember 2020 pukul 08.00-11.00 WIB. Dilanjutkan Majapahit The Birth Of Nusantara Kelana Budaya pukul 21.30 WIB`n#SemangatBaruMarta`n#BelajarDariRumah`n#MediaPemersatuBangsa`n#TVRI #TVRINASIONAL https://t.co/hfmYP0IF44"
$ws.Range("C265").Value = "TheMartaSaputra"
$ws.Range("D265").Value = "Wed Dec 23 22:11:09 +0000 2020"
$ws.Range("A266").Value = [double]"1.3417836187843789E+18"
$ws.Range("B266").Value = "Jadwal Program ""BELAJAR DARI RUMAH""  Hari Kamis, 24 Desember 2020 pukul 08.00-11.00 WIB. Dilanjutkan Majapahit The Birth Of Nusantara Kelana Budaya pukul 21.30 WIB`n#BelajarDariRumah`n#MediaPemersatuBangsa`n#TVRI #TVRINASIONAL https://t.co/oFiwsMyk8S"
$ws.Range("C266").Value = "TVRINasional"
$ws.Range("D266").Value = "Wed Dec 23 16:31:58 +0000 2020"
$ws.Range("A267").Value = [double]"1.34172235418386E+18"
$ws.Range("B267").Value = "Dalam bahasa Inggris, #GenPrestasi biasanya akan mempelajari pronouns dalam grammar. Apa itu Pronouns?`n#BelajarBarengIndiHomeStudy #IndiHomeStudyByIndiHome #BelajarLebihMudah #BelajarTanpaBatas #BelajarDariRumah #dirumahaja #BahasaInggrisSMA #Pronouns https://t.co/1rsqfJ8MKk"
$ws.Range("C267").Value = "indihome_study"
$ws.Range("D267").Value = "Wed Dec 23 12:28:31 +0000 2020"
$ws.Range("A268").Value = [double]"1.3415351941095511E+18"
$ws.Range("B268").Value = "Jadwal Program ""BELAJAR DARI RUMAH""  Hari Rabu, 23 Desember 2020 pukul 08.00-11.00 WIB. Dilanjutkan Panggung Akhir Sekolah Ragam Indonesia pukul 21.30 WIB`n#SemangatBaruMarta`n#BelajarDariRumah`n#MediaPemersatuBangsa`n#TVRI #TVRINASIONAL https://t.co/vwX7YjBbkU"
$ws.Range("C268").Value = "TheMartaSaputra"
$ws.Range("D268").Value = "Wed Dec 23 00:04:49 +0000 2020"
$ws.Range("A269").Value = [double]"1.341526453695906E+18"
$ws.Range("B269").Value = "Selamat pagi, #SahabatDikbud. Jangan lupa sarapan untuk memulai hari, ya! Tayangan-tayangan #BelajardariRumah di @TVRINasional siap hadir untuk menemani waktu belajar #SahabatDikbud. Yuk, simak jadwalnya! `n#MerdekaBelajar`n#BersamaHadapiKorona https://t.co/gR0uWajJXf"
$ws.Range("C269").Value = "Kemdikbud_RI"
$ws.Range("D269").Value = "Tue Dec 22 23:30:05 +0000 2020"
$ws.Range("A270").Value = [double]"1.3414131873510359E+18"
$ws.Range("B270").Value = "Jadwal Program ""BELAJAR DARI RUMAH""  Hari Rabu, 23 Desember 2020 pukul 08.00-11.00 WIB. Dilanjutkan Panggung Akhir Sekolah Ragam Indonesia pukul 21.30 WIB`n#BelajarDariRumah`n#MediaPemersatuBangsa`n#TVRI #TVRINASIONAL https://t.co/t09kry5tqt"
$ws.Range("C270").Value = "TVRINasional"
$ws.Range("D270").Value = "Tue Dec 22 16:00:00 +0000 2020"
$ws.Range("A271").Value = [double]"1.3412255743757271E+18"
$ws.Range("B271").Value = "Baca postingannya di web blog Catatan IzRuYan`nhttps://t.co/n1Smew0E6I`n#IzRuWeb #CatatanIzRuYan #SuaraHatiAnak #BelajardiRumah #BelajardariRumah #Sedih #Susah #Sulit #Cerita #CeritaSedih #CatatanSedih #Curhat #CurhatanAnakSekolah https://t.co/TinaF9ZkPg"
$ws.Range("C271").Value = "IzRuWeb"
$ws.Range("D271").Value = "Tue Dec 22 03:34:30 +0000 2020"
$ws.Range("A272").Value = [double]"1.341197160319169E+18"
$ws.Range("B272").Value = "Terima kasih untuk setiap pelukan, motivasi, dan cinta yang telah kau berikan padaku. `nSelamat Hari Ibu`nmadaniaschool #earlyyears #ibworldschool #belajardarirumah #hariibu #paudmengubahdunia https://t.co/i74qrbl5Ma"
$ws.Range("C272").Value = "tkmadaniayasmin"
$ws.Range("D272").Value = "Tue Dec 22 01:41:35 +0000 2020"
$ws.Range("A273").Value = [double]"1.3411642269586509E+18"
$ws.Range("B273").Value = "Selamat pagi, #SahabatDikbud. Sudah siap belajar kembali hari ini? Jangan lupa siapkan catatanmu, ya. Yuk, simak jadwal acara #BelajardariRumah di @TVRINasional untuk hari Selasa, 22 Desember 2020! `n#MerdekaBelajar`n#BersamaHadapiKorona https://t.co/RcMJb8rIYW"
$ws.Range("C273").Value = "Kemdikbud_RI"
$ws.Range("D273").Value = "Mon Dec 21 23:30:43 +0000 2020"

# --- Split data / sampling-quota duplication: re-append rows 227-233 as 274-280 ---
$ws.Range("A227:D233").Copy($ws.Range("A274"))

# --- Restore view state (selection) ---
$ws.Range("I275").Select()
